$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- original row 15 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D2").Value = 44162
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("Q2").Value = '$/caja 16 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1062
$ws.Range("T2").Value = 16

# Row 3 <- original row 16 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D3").Value = 44162
$ws.Range("K3").Value = 'Castle Brite'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = '$/caja 16 kilos granel'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 938
$ws.Range("T3").Value = 16

# Row 4 <- original row 5 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D4").Value = 44189
$ws.Range("K4").Value = 'Dina'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S4").Value = 833
$ws.Range("T4").Value = 18

# Row 5 <- original row 2 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D5").Value = 44181
$ws.Range("K5").Value = 'Dina'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 220
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 17000
$ws.Range("Q5").Value = '$/caja 18 kilos'
$ws.Range("R5").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S5").Value = 944
$ws.Range("T5").Value = 18

# Row 6 <- original row 17 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D6").Value = 44176
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Provincia de Limarí'
$ws.Range("S6").Value = 944
$ws.Range("T6").Value = 18

# Row 7 <- original row 6 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D7").Value = 44179
$ws.Range("K7").Value = 'Dina'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 18

# Row 8 <- original row 14 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D8").Value = 44160
$ws.Range("K8").Value = 'Dina'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("R8").Value = 'Provincia de Limarí'
$ws.Range("S8").Value = 1333
$ws.Range("T8").Value = 15

# Row 9 <- original row 8 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D9").Value = 44167
$ws.Range("K9").Value = 'Castle Brite'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("Q9").Value = '$/caja 16 kilos granel'
$ws.Range("R9").Value = 'Provincia de Limarí'
$ws.Range("S9").Value = 938
$ws.Range("T9").Value = 16

# Row 10 <- original row 3 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D10").Value = 44174
$ws.Range("K10").Value = 'Castle Brite'
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = '$/caja 10 kilos'
$ws.Range("R10").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S10").Value = 1500
$ws.Range("T10").Value = 10

# Row 11 <- original row 12 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D11").Value = 44168
$ws.Range("K11").Value = 'Castle Brite'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("Q11").Value = '$/caja 10 kilos'
$ws.Range("R11").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 10

# Row 12 <- original row 13 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D12").Value = 44168
$ws.Range("K12").Value = 'Castle Brite'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 17000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Provincia de Limarí'
$ws.Range("S12").Value = 944
$ws.Range("T12").Value = 18

# Row 13 <- original row 7 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D13").Value = 44161
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 1111
$ws.Range("T13").Value = 18

# Row 14 <- original row 4 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D14").Value = 44186
$ws.Range("K14").Value = 'Dina'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("Q14").Value = '$/caja 18 kilos'
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 833
$ws.Range("T14").Value = 18

# Row 15 <- original row 9 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D15").Value = 44172
$ws.Range("K15").Value = 'Castle Brite'
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("Q15").Value = '$/caja 10 kilos'
$ws.Range("R15").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S15").Value = 1500
$ws.Range("T15").Value = 10

# Row 16 <- original row 10 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D16").Value = 44172
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 250
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11000
$ws.Range("P16").Value = 11000
$ws.Range("Q16").Value = '$/caja 10 kilos'
$ws.Range("R16").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S16").Value = 1100
$ws.Range("T16").Value = 10

# Row 17 <- original row 11 data (D,K,L,M,N,O,P,Q,R,S,T)
$ws.Range("D17").Value = 44187
$ws.Range("K17").Value = 'Dina'
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 120
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 889
$ws.Range("T17").Value = 18

